# Updates the cryptos worksheet with the latest scraped price/
# volume snapshot (GitHub Actions scheduled refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column holds text that looks numeric (dot-separated
# thousands, e.g. "27.200.07" / "313.40"). Force it to Text format
# first so Excel stores the values verbatim instead of coercing them
# into numbers (which would drop trailing zeros / collapse the dots).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '27.200.07'
$ws.Range("E2").Value = '  +0.13%  '
$ws.Range("D3").Value = '1.770.46'
$ws.Range("E3").Value = '  +3.42%  '
$ws.Range("D4").Value = '1.003'
$ws.Range("E4").Value = '  +0.25%  '
$ws.Range("D5").Value = '313.40'
$ws.Range("E5").Value = '  +1.67%  '
$ws.Range("E6").Value = '  +0.09%  '
$ws.Range("D7").Value = '0.5260'
$ws.Range("E7").Value = '  +10.10%  '
$ws.Range("D8").Value = '0.3656'
$ws.Range("E8").Value = '  +6.37%  '
$ws.Range("D9").Value = '42.83'
$ws.Range("E9").Value = '  +1.93%  '
$ws.Range("D10").Value = '0.07342'
$ws.Range("E10").Value = '  +0.76%  '
$ws.Range("D11").Value = '1.086'
$ws.Range("E11").Value = '  +3.65%  '
$ws.Range("D12").Value = '1.003'
$ws.Range("E12").Value = '  +0.35%  '
$ws.Range("D13").Value = '20.33'
$ws.Range("E13").Value = '  +2.36%  '
$ws.Range("D14").Value = '6.044'
$ws.Range("E14").Value = '  +3.26%  '
$ws.Range("D15").Value = '1.768.71'
$ws.Range("E15").Value = '  +3.32%  '
$ws.Range("D16").Value = '6.922'
$ws.Range("E16").Value = '  +1.04%  '
$ws.Range("D17").Value = '88.59'
$ws.Range("E17").Value = '  -0.55%  '
$ws.Range("D18").Value = '0.00001041'
$ws.Range("E18").Value = '  +0.29%  '
$ws.Range("D19").Value = '0.06418'
$ws.Range("E19").Value = '  +1.10%  '
$ws.Range("D20").Value = '1.002'
$ws.Range("E20").Value = '  +0.17%  '
$ws.Range("D21").Value = '16.66'
$ws.Range("E21").Value = '  +1.15%  '
$ws.Range("D22").Value = '5.808'
$ws.Range("E22").Value = '  +3.93%  '
$ws.Range("D23").Value = '27.312.96'
$ws.Range("E23").Value = '  +0.46%  '
$ws.Range("E24").Value = '  +3.59%  '
$ws.Range("D25").Value = '2.107'
$ws.Range("E25").Value = '  +0.59%  '
$ws.Range("D26").Value = '154.23'
$ws.Range("E26").Value = '  -0.05%  '
$ws.Range("D27").Value = '20.05'
$ws.Range("E27").Value = '  +1.71%  '
$ws.Range("B28").Value = 'LidoDAOToken'
$ws.Range("C28").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D28").Value = '2.327'
$ws.Range("E28").Value = '  +11.67%  '
$ws.Range("B29").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C29").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D29").Value = '1.971.16'
$ws.Range("E29").Value = '  +3.43%  '
$ws.Range("D30").Value = '120.67'
$ws.Range("E30").Value = '  +0.81%  '
$ws.Range("D31").Value = '1.057'
$ws.Range("E31").Value = '  +4.42%  '
$ws.Range("D32").Value = '0.09711'
$ws.Range("E32").Value = '  +4.47%  '
$ws.Range("D33").Value = '5.541'
$ws.Range("E33").Value = '  +4.39%  '
$ws.Range("D34").Value = '3.623'
$ws.Range("E34").Value = '  +1.16%  '
$ws.Range("D35").Value = '0.02217'
$ws.Range("E35").Value = '  +1.03%  '
$ws.Range("D36").Value = '0.05942'
$ws.Range("E36").Value = '  +1.18%  '
$ws.Range("D37").Value = '11.15'
$ws.Range("E37").Value = '  +1.01%  '
$ws.Range("D38").Value = '4.820'
$ws.Range("E38").Value = '  +1.72%  '
$ws.Range("D39").Value = '0.2013'
$ws.Range("E39").Value = '  +0.90%  '
$ws.Range("D40").Value = '0.6110'
$ws.Range("E40").Value = '  +3.50%  '
$ws.Range("D41").Value = '1.433'
$ws.Range("E41").Value = '  +1.95%  '
$ws.Range("D42").Value = '8.038'
$ws.Range("E42").Value = '  +7.89%  '
$ws.Range("D43").Value = '1.133'
$ws.Range("E43").Value = '  +2.22%  '
$ws.Range("D44").Value = '13.11'
$ws.Range("E44").Value = '  +3.62%  '
$ws.Range("D45").Value = '3.622'
$ws.Range("E45").Value = '  +1.62%  '
$ws.Range("D46").Value = '0.5722'
$ws.Range("E46").Value = '  +1.78%  '
$ws.Range("D47").Value = '120.60'
$ws.Range("E47").Value = '  +2.25%  '
$ws.Range("D48").Value = '1.880'
$ws.Range("E48").Value = '  +2.09%  '
$ws.Range("B49").Value = 'Cronos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D49").Value = '0.06714'
$ws.Range("E49").Value = '  +1.35%  '
$ws.Range("B50").Value = 'EOS'
$ws.Range("C50").Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range("D50").Value = '1.107'
$ws.Range("E50").Value = '  +1.92%  '
$ws.Range("B51").Value = 'PaxDollar'
$ws.Range("C51").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D51").Value = '1.001'
$ws.Range("E51").Value = '  +0.15%  '
